# "minor change on figure"
# Shifts the 15 diagram shapes on slide 4 by a uniform offset
# (dx = +711200 EMU = +56 pt, dy = -13546 EMU = -1.0666... pt).
#
# NOTE: in this runtime, Shape.Left / Shape.Top write straight through to
# the raw <a:off> (they are not rotation-compensated on write, even though
# the getter reports a rotated bounding box for rotated shapes). The host
# also stores Left/Top as f32 and converts to EMU by truncating
# f32(points) * 12700 toward zero, so each literal below is the midpoint
# of the narrow points-range that truncates to the exact target EMU offset
# (the naive target_emu/12700 value can land 1 EMU off after the f32
# round-trip + truncation, especially for negative offsets).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$targets = @{
    2  = @{ Left = 411.0224761948819; Top = 223.43421172834644 }   # Rectangle 1
    3  = @{ Left = 91.99901577795276; Top = 223.43421172834644 }   # Rectangle 2
    4  = @{ Left = 255.0079956059055; Top = 17.633897817716537 }   # Rectangle 3
    6  = @{ Left = 383.1547699094488; Top = 349.37011719015743 }   # Connector: Elbow 16
    7  = @{ Left = 323.53059388110233; Top = -62.60744097480315 }   # Right Brace 6
    10 = @{ Left = 91.99901577795276; Top = 123.32295228582677 }   # Rectangle 9
    11 = @{ Left = 411.02255248503934; Top = 123.32295228582677 }   # Rectangle 10
    16 = @{ Left = 470.937911965748; Top = 464.34901427795273 }   # Rectangle 15
    23 = @{ Left = 411.0224761948819; Top = 332.6279907059055 }   # Rectangle 22
    24 = @{ Left = 255.0079956059055; Top = 464.34901427795273 }   # Rectangle 23
    36 = @{ Left = 179.70027162047245; Top = 174.9552307503937 }   # Straight Arrow Connector 35
    39 = @{ Left = 498.13357545708664; Top = 176.0173645346457 }   # Straight Arrow Connector 38
    48 = @{ Left = 169.34114076220476; Top = 291.5710601720472 }   # Connector: Elbow 16
    57 = @{ Left = 498.13357545708664; Top = 281.212005603937 }   # Straight Arrow Connector 56
    20 = @{ Left = 558.0490112279526; Top = 361.51689149370077 }   # Connector: Elbow 16
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    $id = $shp.Id
    if ($targets.ContainsKey($id)) {
        $t = $targets[$id]
        $shp.Left = $t.Left
        $shp.Top = $t.Top
    }
}
